$wb = $excel.ActiveWorkbook

# --- CmsWork sheet ---
$wsWork = $wb.Worksheets.Item("CmsWork")

# Row 2 (CmsCollection0CmsWork1)
$wsWork.Range("C2").Value = "_:Ne8227b46c9974d46a3bcc9f6ae3abd14"
$wsWork.Range("G2").Value = "CmsCollection0CmsWork1Id1"

# Row 3 (CmsCollection0CmsWork3)
$wsWork.Range("C3").Value = "_:N6c85af05453c4df69069f1f97884180a"
$wsWork.Range("D3").Value = "http://example.com/person0"
$wsWork.Range("U3").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:104"

# Row 4 (CmsCollection1CmsWork5)
$wsWork.Range("C4").Value = "_:Nbca42478582846a885adda29847afff4"
$wsWork.Range("J4").Value = "CmsCollection1CmsWork5 provenance 0"

# Row 5 (CmsCollection1CmsWork7)
$wsWork.Range("C5").Value = "_:N7a84f1163ffe4709a681ce5e742ec193"
$wsWork.Range("D5").Value = "http://example.com/person4"
$wsWork.Range("E5").Value = "CmsCollection1CmsWork7 alternative title 1"
$wsWork.Range("G5").Value = "CmsCollection1CmsWork7Id1"

# Row 6 (FreestandingWork9)
$wsWork.Range("B6").Value = "_:N5c78c639b9714b638b8b2b662c6e0563"
$wsWork.Range("D6").Value = "FreestandingWork9 alternative title 0"
$wsWork.Range("F6").Value = "FreestandingWork9Id0"

# Row 7 (FreestandingWork11)
$wsWork.Range("B7").Value = "_:Nd29b075073c040fb800a676c6295b76a"
$wsWork.Range("C7").Value = "http://example.com/organization3"
$wsWork.Range("T7").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:101"

# --- CmsWorkClosing sheet ---
$wsClosing = $wb.Worksheets.Item("CmsWorkClosing")

$wsClosing.Range("A2").Value = "_:Ndd1dc3ce897545999378fdb08c840e4a"
$wsClosing.Range("C2").Value = "_:N0dd9e86983b445b4a68134ec49e4c210"

$wsClosing.Range("A3").Value = "_:N0f84c8205f394a61819703b8620485c3"
$wsClosing.Range("C3").Value = "_:Nd2382d3861df4467a336dcc71d49d902"

$wsClosing.Range("A4").Value = "_:Neaf5fae07dae4f91a468e2cfa081a066"
$wsClosing.Range("C4").Value = "_:Nf2689075158740228e20d07fc423c200"

$wsClosing.Range("A5").Value = "_:N24af2fc7193a419fa584727384fde638"
$wsClosing.Range("C5").Value = "_:Nf3dcefad2d19487a8f3b687f4059735b"

$wsClosing.Range("A6").Value = "_:N49c6d7509f6a43beb36ac7ca6904a2f9"
$wsClosing.Range("C6").Value = "_:Nb9eb50f11645437e99fc0a3081786a39"

$wsClosing.Range("A7").Value = "_:Ndd7a0563da6841de827d3ae55d053a32"
$wsClosing.Range("C7").Value = "_:Nab77dd39560c403491720856b0469177"

# --- CmsWorkOpening sheet ---
$wsOpening = $wb.Worksheets.Item("CmsWorkOpening")

$wsOpening.Range("C2").Value = "_:N0dd9e86983b445b4a68134ec49e4c210"
$wsOpening.Range("C3").Value = "_:Nd2382d3861df4467a336dcc71d49d902"
$wsOpening.Range("C4").Value = "_:Nf2689075158740228e20d07fc423c200"
$wsOpening.Range("C5").Value = "_:Nf3dcefad2d19487a8f3b687f4059735b"
$wsOpening.Range("C6").Value = "_:Nb9eb50f11645437e99fc0a3081786a39"
$wsOpening.Range("C7").Value = "_:Nab77dd39560c403491720856b0469177"

# --- CmsRightsStatement sheet ---
$wsRights = $wb.Worksheets.Item("CmsRightsStatement")
$wsRights.Range("E2").Value = "You may find additional information about the copyright status of the Item on the website of the organization that has made the Item available."
